# HW rev. 1 specified BOM, ERC & DRC check out, generated gerbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old data rows (2..9); header row 1 is kept
$ws.Rows("2:9").Delete()

# Add new header columns for mouser part no. / mouser url
$ws.Range("F1").Value = "mouser part no."
$ws.Range("G1").Value = "mouser url"

# New BOM rows (rows 2-8): Qty, Value, Device, Package, Parts
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "0.1uF"
$ws.Cells.Item(2, 3).Value = "C-EUC0402"
$ws.Cells.Item(2, 4).Value = "C0402"
$ws.Cells.Item(2, 5).Value = "C8"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "10k"
$ws.Cells.Item(3, 3).Value = "R-EU_R0402"
$ws.Cells.Item(3, 4).Value = "R0402"
$ws.Cells.Item(3, 5).Value = "R1, R2"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "10uF"
$ws.Cells.Item(4, 3).Value = "C-EUC0603"
$ws.Cells.Item(4, 4).Value = "C0603"
$ws.Cells.Item(4, 5).Value = "C3, C4"

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "10uF"
$ws.Cells.Item(5, 3).Value = "C-EUC0805"
$ws.Cells.Item(5, 4).Value = "C0805"
$ws.Cells.Item(5, 5).Value = "C1, C2"

$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = "22uF"
$ws.Cells.Item(6, 3).Value = "C-EUC0805"
$ws.Cells.Item(6, 4).Value = "C0805"
$ws.Cells.Item(6, 5).Value = "C5, C6, C7"

$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "TPS630701RNMT"
$ws.Cells.Item(7, 3).Value = "TPS630701"
$ws.Cells.Item(7, 4).Value = "VQFN"
$ws.Cells.Item(7, 5).Value = "IC1"

$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "XFL4020-152MEC  - 1.5uH"
$ws.Cells.Item(8, 3).Value = "XFL4020-152MEC"
$ws.Cells.Item(8, 4).Value = "SMD"
$ws.Cells.Item(8, 5).Value = "L1"

# mouser part no. / mouser url columns, filled row-by-row in the order
# the parts were looked up (IC1, L1, C5-C7, C1-C2, C3-C4, C8, R1-R2)
$ws.Cells.Item(7, 6).Value = "595-TPS630701RNMR "
$ws.Hyperlinks.Add($ws.Cells.Item(7, 7), "https://hr.mouser.com/ProductDetail/Texas-Instruments/TPS630701RNMR?qs=LuYMPh7GGMTyE4ON9PhzXw%3D%3D")

$ws.Cells.Item(8, 6).Value = "994-XFL4020-152MEC "
$ws.Hyperlinks.Add($ws.Cells.Item(8, 7), "https://hr.mouser.com/ProductDetail/Coilcraft/XFL4020-152MEC?qs=%2Fha2pyFaduigLhcKZCbSe9QbArGlvPlxXOCFjm%2FaTBqm0TUkZLvsng%3D%3D")

$ws.Cells.Item(6, 6).Value = "80-C0805C226M8"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 7), "https://hr.mouser.com/ProductDetail/KEMET/C0805C226M8PACTU?qs=cGEy3R83DS93ZizaMBlKFQ%3D%3D")

$ws.Cells.Item(5, 6).Value = "80-C0805C106M3PACLR "
$ws.Hyperlinks.Add($ws.Cells.Item(5, 7), "https://hr.mouser.com/ProductDetail/KEMET/C0805C106M3PAC7210?qs=sSYV1F9c5cFaojLA0ITCpw%3D%3D")

$ws.Cells.Item(4, 6).Value = "81-GRM188R61E106MA3J"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 7), "https://hr.mouser.com/ProductDetail/Murata-Electronics/GRM188R61E106MA73J?qs=hNud%2FORuBR2%252B%252BY67hhW1nA%3D%3D")

$ws.Cells.Item(2, 6).Value = "80-C0402C104K3P "
$ws.Hyperlinks.Add($ws.Cells.Item(2, 7), "https://hr.mouser.com/ProductDetail/KEMET/C0402C104K3PACTU?qs=gt1LBUVyoHnkmt1KfrmtmQ%3D%3D")

$ws.Cells.Item(3, 6).Value = "603-RC0402JR-7W10KL "
$ws.Hyperlinks.Add($ws.Cells.Item(3, 7), "https://hr.mouser.com/ProductDetail/Yageo/RC0402JR-7W10KL?qs=sGAEpiMZZMvdGkrng054t3bQBroXiaAZcfok9fjjif4OfQngMUmcIg%3D%3D")

# Column G is now much wider to fit the URLs (107.1 round-trips to a
# stored width of exactly 108 in the saved XML)
$ws.Columns("G").ColumnWidth = 107.1

# Restore selection to match post-edit state
$ws.Range("F14").Select()

$wb.Save()
